$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.294.71'
$ws.Range("E2").Value = '  +1.21%  '
$ws.Range("D3").Value = '1.620.63'
$ws.Range("E3").Value = '  +1.67%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.16'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.72%  '
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("E7").Value = '  +1.04%  '
$ws.Range("E8").Value = '  +0.96%  '
$ws.Range("E11").Value = '  +0.73%  '
$ws.Range("D12").Value = '1.846.17'
$ws.Range("E12").Value = '  +1.66%  '
$ws.Range("D13").Value = '1.624.76'
$ws.Range("E13").Value = '  +1.78%  '
$ws.Range("E14").Value = '  +0.19%  '
$ws.Range("E15").Value = '  +1.18%  '
$ws.Range("D16").Value = '26.305.56'
$ws.Range("E16").Value = '  +1.22%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.38'
$ws.Range("D17").Style = "Normal"
$ws.Range("E18").Value = '  +0.93%  '
$ws.Range("E19").Value = '  +0.01%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '201.73'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.93%  '
$ws.Range("E21").Value = '  +1.23%  '
$ws.Range("E22").Value = '  +1.19%  '
$ws.Range("E23").Value = '  +0.82%  '
$ws.Range("E24").Value = '  +4.44%  '
$ws.Range("E25").Value = '  +0.58%  '
$ws.Range("E26").Value = '  +0.08%  '
$ws.Range("E27").Value = '  -0.41%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.17'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.68%  '
$ws.Range("E29").Value = '  +1.69%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0526'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +10.81%  '
$ws.Range("E31").Value = '  +0.62%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.17'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.44%  '
$ws.Range("E33").Value = '  +0.04%  '
$ws.Range("E34").Value = '  +1.20%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.42'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.33%  '
$ws.Range("D36").Value = '1.180.59'
$ws.Range("E36").Value = '  +5.13%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0164'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.03%  '
$ws.Range("E38").Value = '  +3.26%  '
$ws.Range("E39").Value = '  +0.06%  '
$ws.Range("E40").Value = '  +0.20%  '
$ws.Range("E41").Value = '  +0.98%  '
$ws.Range("E42").Value = '  +0.57%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.34'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.83%  '
$ws.Range("D44").Value = '1.757.93'
$ws.Range("E44").Value = '  +1.78%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '93.28'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.78%  '
$ws.Range("E46").Value = '  +15.13%  '
$ws.Range("E47").Value = '  +1.94%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '53.85'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.99%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0508'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.06%  '
$ws.Range("E50").Value = '  +0.31%  '
$ws.Range("E51").Value = '  -0.26%  '
